$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Weekly Quantity": append two new weeks at the bottom (rows 49-50)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$dateFmtWeekly = $wsWeekly.Cells.Item(2, 1).NumberFormat

$wsWeekly.Cells.Item(49, 1).Value = 45662.99999999999
$wsWeekly.Cells.Item(49, 1).NumberFormat = $dateFmtWeekly
$wsWeekly.Cells.Item(49, 2).Value = 6

$wsWeekly.Cells.Item(50, 1).Value = 45669.99999999999
$wsWeekly.Cells.Item(50, 1).NumberFormat = $dateFmtWeekly
$wsWeekly.Cells.Item(50, 2).Value = 1

# ---------------------------------------------------------------------
# Sheet "Monthly Trend": append one new month at the bottom (row 21)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$dateFmtMonthly = $wsMonthly.Cells.Item(2, 1).NumberFormat

$wsMonthly.Cells.Item(21, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(21, 1).NumberFormat = $dateFmtMonthly
$wsMonthly.Cells.Item(21, 2).Value = 7

# ---------------------------------------------------------------------
# Sheet "PO Forecast": refreshed forecast model
#   - several existing forecast values shift down by one
#   - the forecast's "flat tail" (constant value 8) slides forward in time
#     and gains two extra weeks at the end
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")
$dateFmtForecast = $wsForecast.Cells.Item(2, 1).NumberFormat

# Updated forecast quantities (dates unchanged) for rows 12-47
$qtyUpdates = @{
    12 = 5
    13 = 5
    21 = 5
    29 = 6
    30 = 6
    31 = 6
    43 = 7
    44 = 7
    45 = 7
    46 = 7
    47 = 7
}
foreach ($row in $qtyUpdates.Keys) {
    $wsForecast.Cells.Item($row, 2).Value = $qtyUpdates[$row]
}

# Rows 49-58: the trailing flat-8 forecast shifts forward by two weeks and
# extends two rows further (row 48 at 45641.99999999999 stays as-is).
$tailDates = @{
    49 = 45662.99999999999
    50 = 45669.99999999999
    51 = 45676.99999999999
    52 = 45683.99999999999
    53 = 45690.99999999999
    54 = 45697.99999999999
    55 = 45704.99999999999
    56 = 45711.99999999999
    57 = 45718.99999999999
    58 = 45725.99999999999
}
foreach ($row in (49..58)) {
    $wsForecast.Cells.Item($row, 1).Value = $tailDates[$row]
    $wsForecast.Cells.Item($row, 1).NumberFormat = $dateFmtForecast
    $wsForecast.Cells.Item($row, 2).Value = 8
}
